$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4545454545454545
$ws.Range("D2").Value = 0.4761904761904762
$ws.Range("B3").Value = 0.6153846153846154
$ws.Range("C3").Value = 0.5714285714285714
$ws.Range("D3").Value = 0.5925925925925927
$ws.Range("B4").Value = 0.5416666666666666
$ws.Range("C4").Value = 0.5416666666666666
$ws.Range("D4").Value = 0.5416666666666666
$ws.Range("E4").Value = 0.5416666666666666
$ws.Range("B5").Value = 0.534965034965035
$ws.Range("C5").Value = 0.5357142857142857
$ws.Range("D5").Value = 0.5343915343915344
$ws.Range("B6").Value = 0.5483682983682984
$ws.Range("C6").Value = 0.5416666666666666
$ws.Range("D6").Value = 0.5440917107583775
$ws.Range("B7").Value = 0.4545454545454545
$ws.Range("D7").Value = 0.4761904761904762
$ws.Range("B8").Value = 0.6153846153846154
$ws.Range("C8").Value = 0.5714285714285714
$ws.Range("D8").Value = 0.5925925925925927
$ws.Range("B9").Value = 0.5416666666666666
$ws.Range("C9").Value = 0.5416666666666666
$ws.Range("D9").Value = 0.5416666666666666
$ws.Range("E9").Value = 0.5416666666666666
$ws.Range("B10").Value = 0.534965034965035
$ws.Range("C10").Value = 0.5357142857142857
$ws.Range("D10").Value = 0.5343915343915344
$ws.Range("B11").Value = 0.5483682983682984
$ws.Range("C11").Value = 0.5416666666666666
$ws.Range("D11").Value = 0.5440917107583775
$ws.Range("B12").Value = 0.6
$ws.Range("C12").Value = 0.6
$ws.Range("D12").Value = 0.6
$ws.Range("B13").Value = 0.7142857142857143
$ws.Range("D13").Value = 0.7142857142857143
$ws.Range("B14").Value = 0.6666666666666666
$ws.Range("C14").Value = 0.6666666666666666
$ws.Range("D14").Value = 0.6666666666666666
$ws.Range("E14").Value = 0.6666666666666666
$ws.Range("B15").Value = 0.6571428571428571
$ws.Range("C15").Value = 0.6571428571428571
$ws.Range("D15").Value = 0.6571428571428571
$ws.Range("B16").Value = 0.6666666666666666
$ws.Range("C16").Value = 0.6666666666666666
$ws.Range("D16").Value = 0.6666666666666666
$ws.Range("B17").Value = 0.3846153846153846
$ws.Range("D17").Value = 0.4347826086956522
$ws.Range("B18").Value = 0.5454545454545454
$ws.Range("C18").Value = 0.4285714285714285
$ws.Range("D18").Value = 0.4799999999999999
$ws.Range("B19").Value = 0.4583333333333333
$ws.Range("C19").Value = 0.4583333333333333
$ws.Range("D19").Value = 0.4583333333333333
$ws.Range("E19").Value = 0.4583333333333333
$ws.Range("B20").Value = 0.465034965034965
$ws.Range("C20").Value = 0.4642857142857143
$ws.Range("D20").Value = 0.457391304347826
$ws.Range("B21").Value = 0.4784382284382284
$ws.Range("C21").Value = 0.4583333333333333
$ws.Range("D21").Value = 0.4611594202898551
$ws.Range("B22").Value = 0.5454545454545454
$ws.Range("D22").Value = 0.5714285714285713
$ws.Range("B23").Value = 0.6923076923076923
$ws.Range("C23").Value = 0.6428571428571429
$ws.Range("D23").Value = 0.6666666666666666
$ws.Range("B24").Value = 0.625
$ws.Range("C24").Value = 0.625
$ws.Range("D24").Value = 0.625
$ws.Range("E24").Value = 0.625
$ws.Range("B25").Value = 0.6188811188811189
$ws.Range("C25").Value = 0.6214285714285714
$ws.Range("D25").Value = 0.619047619047619
$ws.Range("B26").Value = 0.6311188811188811
$ws.Range("C26").Value = 0.625
$ws.Range("D26").Value = 0.6269841269841269
